$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the train schedule values in row 2 (D2, F2, H2)
$ws.Range("D2").Value = 6
$ws.Range("F2").Value = -3
$ws.Range("H2").Value = 46

# Move the active selection from E5 to C2
[void]$ws.Range("C2").Select()
